# "#55: fixed import and timeslot management"
# Add a header row (Raum / Kapazitaet) above the existing room list, and
# add a new "Kapazitaet" (capacity) column with a capacity value for each room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; everything below shifts down one row.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Raum"
$ws.Range("B1").Value = "Kapazität"

# Re-assert the room number in A2 as text so the leading zero survives
# ("008" must stay text, not become the number 8).
$ws.Range("A2").Value = "'008"

# Capacity values for each room, row by row (rows 2-15 after the insert).
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 20
$ws.Range("B7").Value = 20
$ws.Range("B8").Value = 20
$ws.Range("B9").Value = 20
$ws.Range("B10").Value = 20
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 20
$ws.Range("B14").Value = 50
$ws.Range("B15").Value = 20

# Leave the cursor on the next empty row, as in the edited workbook.
$ws.Rows.Item(16).Select()
